# Update cryptocurrency Price (D) and Volume(1h) (E) columns per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.894.11'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = '1.631.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'211.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.58%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'23.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '1.861.82'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = '1.631.95'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'65.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '27.886.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'229.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = '  +1.34%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = '  -4.45%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'154.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.82%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'6.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '1.389.08'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '  +10.60%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = '  +1.27%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.558"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.94%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = '  -3.38%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'65.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.13%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'5.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '1.772.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = '  -2.79%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'88.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D51").Value = "'7.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.93%  '
$ws.Range("E51").Style = "Normal"
